# The edit targets "Sheet2" (the Video Recording Progress sheet), which is
# already the ActiveSheet/ActiveWorkbook's selected/active tab.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark the remaining progress-tracking checkboxes (columns D:H, rows 479-513)
# as completed (TRUE), matching the diff which flips every boolean cell in
# that block from FALSE (0) to TRUE (1).
$ws.Range("D479:H513").Value = $true

# Update the view state to mirror the diff: the window is scrolled down two
# more rows (topLeftCell A479 -> A481) and the selection now spans the whole
# freshly-completed block C479:H513 (still anchored at C479).
$excel.ActiveWindow.ScrollRow = 481
$ws.Range("C479:H513").Select()
